# ToDoList.xlsx small-defects update, per commit "Small defecs as per todo list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status updates: several items flip from "Working" to "Done" ---
$ws.Range("C6").Value = "Done"
$ws.Range("C8").Value = "Done"
$ws.Range("C9").Value = "Done"
$ws.Range("C12").Value = "Done"

# --- New todo row 14 (Sr.Number 13) ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Please change the backgroud color."
$ws.Range("B14").WrapText = $true
$ws.Range("C14").Value = "Working"
$ws.Range("D14").Value = "Sharmila"

# --- New comment text on row 12, wrapped like the other note cells ---
$ws.Range("E12").Value = "Once session timeout it gets log off. I update the session timeout after 30min from ideal. If its still not working please let me know."
$ws.Range("E12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 60

# --- New "Comment" column header (E1) ---
$ws.Range("E1").Value = "Comment"
$ws.Columns.Item(5).ColumnWidth = 32.5

# --- New todo row 15 (Sr.Number 14) ---
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Default save message should get removed."
$ws.Range("B15").WrapText = $true
$ws.Range("C15").Value = "Done"
$ws.Range("D15").Value = "Taniya"

# --- New todo row 16 (Sr.Number 15) ---
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Submit button and Button style should same"
$ws.Range("B16").WrapText = $true
$ws.Range("C16").Value = "Done"
$ws.Range("D16").Value = "Taniya"

# --- Restore view: scroll to top, select D1 ---
$ws.Range("D1").Select()
